# Updates cryptocurrency price/volume table (columns D and E) on Sheet1
# to reflect the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.202.49"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.904.89"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'306.74"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.5253"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("D8").Value = "'0.3779"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "'0.07258"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").Value = "'21.14"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "'0.8998"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "'0.08427"
$ws.Range("E12").Value = "  +10.55%  "
$ws.Range("D13").Value = "1.892.66"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'94.73"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'5.270"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "'0.000008628"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "'14.57"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D20").Value = "27.235.66"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "'5.060"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "2.137.19"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'6.435"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'146.89"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("D26").Value = "'2.273"
$ws.Range("E26").Value = "  +5.45%  "
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "'114.89"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'4.928"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "'4.806"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'0.09295"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").Value = "'0.8073"
$ws.Range("E33").Value = "  +6.39%  "
$ws.Range("D34").Value = "'0.05064"
$ws.Range("D35").Value = "'1.236"
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("D36").Value = "'2.952"
$ws.Range("E36").Value = "  -2.83%  "
$ws.Range("D37").Value = "'3.360"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").Value = "'2.603"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").Value = "'0.5722"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "'0.01987"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").Value = "'1.073"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "'6.646"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "'8.966"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "'117.76"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").Value = "'0.1516"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").Value = "'0.4847"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'10.14"
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("D49").Value = "'1.616"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("D50").Value = "'37.46"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").Value = "'63.72"
$ws.Range("E51").Value = "  +0.24%  "
